$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 183
$ws.Cells.Item(183, 1).Value = "13-09-2021"
$ws.Cells.Item(183, 2).Value = 1.33
$ws.Cells.Item(183, 3).Value = 0.74
$ws.Cells.Item(183, 4).Value = 0.05
$ws.Cells.Item(183, 5).Value = -0.33
$ws.Cells.Item(183, 6).Value = 2.06
$ws.Cells.Item(183, 7).Value = 3.26
$ws.Cells.Item(183, 8).Value = 2
$ws.Cells.Item(183, 9).Value = 1.89
$ws.Cells.Item(183, 10).Value = 6.99
$ws.Cells.Item(183, 11).Value = 1.68
$ws.Cells.Item(183, 12).Value = 6.15
$ws.Cells.Item(183, 13).Value = 10.96
$ws.Cells.Item(183, 14).Value = 7.42
$ws.Cells.Item(183, 15).Value = 6.99
$ws.Cells.Item(183, 16).Value = 6.43
# Row 184
$ws.Cells.Item(184, 1).Value = "14-09-2021"
$ws.Cells.Item(184, 2).Value = 1.29
$ws.Cells.Item(184, 3).Value = 0.74
$ws.Cells.Item(184, 4).Value = 0.05
$ws.Cells.Item(184, 5).Value = -0.34
$ws.Cells.Item(184, 6).Value = 2.09
$ws.Cells.Item(184, 7).Value = 3.31
$ws.Cells.Item(184, 8).Value = 1.98
$ws.Cells.Item(184, 9).Value = 1.9
$ws.Cells.Item(184, 10).Value = 7.03
$ws.Cells.Item(184, 11).Value = 1.67
$ws.Cells.Item(184, 12).Value = 6.16
$ws.Cells.Item(184, 13).Value = 10.94
$ws.Cells.Item(184, 14).Value = 7.43
$ws.Cells.Item(184, 15).Value = 6.98
$ws.Cells.Item(184, 16).Value = 6.42
# Row 185
$ws.Cells.Item(185, 1).Value = "15-09-2021"
$ws.Cells.Item(185, 2).Value = 1.3
$ws.Cells.Item(185, 3).Value = 0.78
$ws.Cells.Item(185, 4).Value = 0.04
$ws.Cells.Item(185, 5).Value = -0.31
$ws.Cells.Item(185, 6).Value = 2.05
$ws.Cells.Item(185, 7).Value = 3.3
$ws.Cells.Item(185, 8).Value = 1.98
$ws.Cells.Item(185, 9).Value = 1.91
$ws.Cells.Item(185, 10).Value = 7.06
$ws.Cells.Item(185, 11).Value = 1.65
$ws.Cells.Item(185, 12).Value = 6.2
$ws.Cells.Item(185, 13).Value = 11.03
$ws.Cells.Item(185, 14).Value = 7.4
$ws.Cells.Item(185, 15).Value = 6.99
$ws.Cells.Item(185, 16).Value = 6.44
# Row 186
$ws.Cells.Item(186, 1).Value = "16-09-2021"
$ws.Cells.Item(186, 2).Value = 1.34
$ws.Cells.Item(186, 3).Value = 0.82
$ws.Cells.Item(186, 4).Value = 0.05
$ws.Cells.Item(186, 5).Value = -0.3
$ws.Cells.Item(186, 6).Value = 2.06
$ws.Cells.Item(186, 8).Value = 2.02
$ws.Cells.Item(186, 9).Value = 1.97
$ws.Cells.Item(186, 10).Value = 7.06
$ws.Cells.Item(186, 11).Value = 1.76
$ws.Cells.Item(186, 12).Value = 6.31
$ws.Cells.Item(186, 13).Value = 11.07
$ws.Cells.Item(186, 14).Value = 7.47
$ws.Cells.Item(186, 16).Value = 6.41
# Row 187
$ws.Cells.Item(187, 1).Value = "17-09-2021"
$ws.Cells.Item(187, 2).Value = 1.36
$ws.Cells.Item(187, 3).Value = 0.85
$ws.Cells.Item(187, 4).Value = 0.05
$ws.Cells.Item(187, 5).Value = -0.28
$ws.Cells.Item(187, 6).Value = 2.1
$ws.Cells.Item(187, 7).Value = 3.32
$ws.Cells.Item(187, 8).Value = 2.04
$ws.Cells.Item(187, 9).Value = 1.99
$ws.Cells.Item(187, 10).Value = 7.09
$ws.Cells.Item(187, 11).Value = 1.8
$ws.Cells.Item(187, 12).Value = 6.41
$ws.Cells.Item(187, 13).Value = 11.14
$ws.Cells.Item(187, 14).Value = 7.53
$ws.Cells.Item(187, 15).Value = 7.04
$ws.Cells.Item(187, 16).Value = 6.38
# Row 188
$ws.Cells.Item(188, 1).Value = "20-09-2021"
$ws.Cells.Item(188, 2).Value = 1.31
$ws.Cells.Item(188, 3).Value = 0.79
$ws.Cells.Item(188, 4).Value = 0.05
$ws.Cells.Item(188, 5).Value = -0.33
$ws.Cells.Item(188, 7).Value = 3.35
$ws.Cells.Item(188, 8).Value = 1.99
$ws.Cells.Item(188, 9).Value = 1.98
$ws.Cells.Item(188, 10).Value = 7.17
$ws.Cells.Item(188, 11).Value = 1.79
$ws.Cells.Item(188, 12).Value = 6.42
$ws.Cells.Item(188, 13).Value = 11.2
$ws.Cells.Item(188, 15).Value = 7.08
$ws.Cells.Item(188, 16).Value = 6.36
